$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.30729222605736
$ws.Range("C2").Value = 7.458549524718691
$ws.Range("D2").Value = 13.87477634477309
$ws.Range("E2").Value = 14.36979348599436
$ws.Range("G2").Value = 3.734287269652385
$ws.Range("J2").Value = 8.680378389478513
$ws.Range("K2").Value = 13.05398466323908
$ws.Range("L2").Value = 12.12290868794631
$ws.Range("M2").Value = 17.74375959756803
$ws.Range("N2").Value = 23.69736278369955
$ws.Range("O2").Value = 34.68587127431238

$ws.Range("B3").Value = 17.18677291564855
$ws.Range("C3").Value = 7.425098116930809
$ws.Range("D3").Value = 13.87777823367869
$ws.Range("E3").Value = 14.39546135602054
$ws.Range("G3").Value = 3.736373415549733
$ws.Range("J3").Value = 8.68551476181595
$ws.Range("K3").Value = 12.96748133574138
$ws.Range("L3").Value = 12.13613344316188
$ws.Range("M3").Value = 17.73866160537809
$ws.Range("N3").Value = 23.7596625166325
$ws.Range("O3").Value = 34.76198623989879

$ws.Range("B4").Value = 17.11581112039163
$ws.Range("C4").Value = 7.404117302737166
$ws.Range("D4").Value = 13.88188508070126
$ws.Range("E4").Value = 14.41272521874182
$ws.Range("G4").Value = 3.737723296768225
$ws.Range("J4").Value = 8.688841695596224
$ws.Range("K4").Value = 12.91651078041265
$ws.Range("L4").Value = 12.14562875516223
$ws.Range("M4").Value = 17.73801024168178
$ws.Range("N4").Value = 23.79977174909328
$ws.Range("O4").Value = 34.81411826974769

$ws.Range("B5").Value = 17.0876838870935
$ws.Range("C5").Value = 7.395458129557345
$ws.Range("D5").Value = 13.8841289223226
$ws.Range("E5").Value = 14.42013901172524
$ws.Range("G5").Value = 3.738290784282964
$ws.Range("J5").Value = 8.690241118769206
$ws.Range("K5").Value = 12.89629699351698
$ws.Range("L5").Value = 12.14984455251038
$ws.Range("M5").Value = 17.7383704322487
$ws.Range("N5").Value = 23.81658477555652
$ws.Range("O5").Value = 34.83671896138645

$ws.Range("B6").Value = 17.08306185973282
$ws.Range("C6").Value = 7.394013727386425
$ws.Range("D6").Value = 13.8845359800412
$ws.Range("E6").Value = 14.42139294957485
$ws.Range("G6").Value = 3.738386067673353
$ws.Range("J6").Value = 8.690476133142814
$ws.Range("K6").Value = 12.89297467620984
$ws.Range("L6").Value = 12.15056551870108
$ws.Range("M6").Value = 17.73846809084901
$ws.Range("N6").Value = 23.81940487989565
$ws.Range("O6").Value = 34.84055369426893

$ws.Range("B7").Value = 17.115428552911
$ws.Range("C7").Value = 7.404000962488285
$ws.Range("D7").Value = 13.88191303169624
$ws.Range("E7").Value = 14.41282367001332
$ws.Range("G7").Value = 3.737730879578614
$ws.Range("J7").Value = 8.688860391721622
$ws.Range("K7").Value = 12.91623589090817
$ws.Range("L7").Value = 12.1456842076803
$ws.Range("M7").Value = 17.73801256332384
$ws.Range("N7").Value = 23.79999659812694
$ws.Range("O7").Value = 34.81441757930079

$ws.Range("B8").Value = 17.26512298280721
$ws.Range("C8").Value = 7.447108132235031
$ws.Range("D8").Value = 13.87534225531797
$ws.Range("E8").Value = 14.37833198883945
$ws.Range("G8").Value = 3.734992289020981
$ws.Range("J8").Value = 8.68211354123695
$ws.Range("K8").Value = 13.02372480395378
$ws.Range("L8").Value = 12.1271834940252
$ws.Range("M8").Value = 17.74148877662572
$ws.Range("N8").Value = 23.71845907509869
$ws.Range("O8").Value = 34.71099515383916

$ws.Range("B9").Value = 17.58155476580446
$ws.Range("C9").Value = 7.528065133270803
$ws.Range("D9").Value = 13.88036536210647
$ws.Range("E9").Value = 14.32260219701752
$ws.Range("G9").Value = 3.730166752569117
$ws.Range("J9").Value = 8.670251482943234
$ws.Range("K9").Value = 13.25067045923755
$ws.Range("L9").Value = 12.10179047319698
$ws.Range("M9").Value = 17.76786514110492
$ws.Range("N9").Value = 23.57324054912487
$ws.Range("O9").Value = 34.55103655370038

$ws.Range("B10").Value = 17.8262794622818
$ws.Range("C10").Value = 7.585269494551017
$ws.Range("D10").Value = 13.89489221211767
$ws.Range("E10").Value = 14.28888619440412
$ws.Range("G10").Value = 3.72695009625276
$ws.Range("J10").Value = 8.662362842355051
$ws.Range("K10").Value = 13.42608217585228
$ws.Range("L10").Value = 12.08973625825656
$ws.Range("M10").Value = 17.79901075355746
$ws.Range("N10").Value = 23.47541349556935
$ws.Range("O10").Value = 34.45967486928514

$ws.Range("B11").Value = 17.93988843721434
$ws.Range("C11").Value = 7.610779312642106
$ws.Range("D11").Value = 13.90383232548671
$ws.Range("E11").Value = 14.27511082100914
$ws.Range("G11").Value = 3.72555738319601
$ws.Range("J11").Value = 8.658951874905609
$ws.Range("K11").Value = 13.50750248109059
$ws.Range("L11").Value = 12.08567790637643
$ws.Range("M11").Value = 17.81569654569219
$ws.Range("N11").Value = 23.43281692365373
$ws.Range("O11").Value = 34.42379808211001

$ws.Range("B12").Value = 17.98320449624672
$ws.Range("C12").Value = 7.620363730786968
$ws.Range("D12").Value = 13.90755079786095
$ws.Range("E12").Value = 14.27011852776014
$ws.Range("G12").Value = 3.725040089378485
$ws.Range("J12").Value = 8.657685645336318
$ws.Range("K12").Value = 13.53854511593351
$ws.Range("L12").Value = 12.08434528484931
$ws.Range("M12").Value = 17.82237335541522
$ws.Range("N12").Value = 23.41695943207879
$ws.Range("O12").Value = 34.41103009517261

$ws.Range("B13").Value = 17.97386303213898
$ws.Range("C13").Value = 7.618302946643772
$ws.Range("D13").Value = 13.90673518438081
$ws.Range("E13").Value = 14.27118374678597
$ws.Range("G13").Value = 3.725151049610312
$ws.Range("J13").Value = 8.657957221290538
$ws.Range("K13").Value = 13.53185053606183
$ws.Range("L13").Value = 12.08462321981757
$ws.Range("M13").Value = 17.82091951344733
$ws.Range("N13").Value = 23.42036250635288
$ws.Range("O13").Value = 34.41374353385309

$ws.Range("B14").Value = 17.94344635978083
$ws.Range("C14").Value = 7.61156935421298
$ws.Range("D14").Value = 13.90413158545828
$ws.Range("E14").Value = 14.27469561281945
$ws.Range("G14").Value = 3.72551462306029
$ws.Range("J14").Value = 8.658847192418357
$ws.Range("K14").Value = 13.51005229116482
$ws.Range("L14").Value = 12.08556418301466
$ws.Range("M14").Value = 17.81623869307485
$ws.Range("N14").Value = 23.43150685506631
$ws.Range("O14").Value = 34.42273126088586

$ws.Range("B15").Value = 17.9248526641826
$ws.Range("C15").Value = 7.60743493430448
$ws.Range("D15").Value = 13.90258010740164
$ws.Range("E15").Value = 14.2768759069489
$ws.Range("G15").Value = 3.725738635666002
$ws.Range("J15").Value = 8.659395633862161
$ws.Range("K15").Value = 13.49672696436717
$ws.Range("L15").Value = 12.08616711848236
$ws.Range("M15").Value = 17.8134180963463
$ws.Range("N15").Value = 23.43836859855036
$ws.Range("O15").Value = 34.42834301126899

$ws.Range("B16").Value = 17.81889793079319
$ws.Range("C16").Value = 7.583591890768486
$ws.Range("D16").Value = 13.89435468776773
$ws.Range("E16").Value = 14.28981784041137
$ws.Range("G16").Value = 3.727042528886284
$ws.Range("J16").Value = 8.662589321649557
$ws.Range("K16").Value = 13.4207918922179
$ws.Range("L16").Value = 12.09003009610658
$ws.Range("M16").Value = 17.79797064118496
$ws.Range("N16").Value = 23.47823553024841
$ws.Range("O16").Value = 34.46213391818622

$ws.Range("B17").Value = 17.75445909100972
$ws.Range("C17").Value = 7.56883230270221
$ws.Range("D17").Value = 13.8899043448403
$ws.Range("E17").Value = 14.29815706249379
$ws.Range("G17").Value = 3.727860461515425
$ws.Range("J17").Value = 8.664593957745568
$ws.Range("K17").Value = 13.37460796217795
$ws.Range("L17").Value = 12.09276440882704
$ws.Range("M17").Value = 17.78913632751701
$ws.Range("N17").Value = 23.5031798485423
$ws.Range("O17").Value = 34.48431971162798

$ws.Range("B18").Value = 17.71761244297026
$ws.Range("C18").Value = 7.560294828406586
$ws.Range("D18").Value = 13.88756425254753
$ws.Range("E18").Value = 14.30310063535364
$ws.Range("G18").Value = 3.728337559140677
$ws.Range("J18").Value = 8.665763695918475
$ws.Range("K18").Value = 13.34819850927591
$ws.Range("L18").Value = 12.09447130462645
$ws.Range("M18").Value = 17.78429226250055
$ws.Range("N18").Value = 23.51770657478886
$ws.Range("O18").Value = 34.49761539166317

$ws.Range("B19").Value = 17.70517505982917
$ws.Range("C19").Value = 7.55739597353324
$ws.Range("D19").Value = 13.88680972395168
$ws.Range("E19").Value = 14.30479971918017
$ws.Range("G19").Value = 3.72850023891493
$ws.Range("J19").Value = 8.666162625014707
$ws.Range("K19").Value = 13.33928392590411
$ws.Range("L19").Value = 12.09507229984546
$ws.Range("M19").Value = 17.78269299471398
$ws.Range("N19").Value = 23.5226559228393
$ws.Range("O19").Value = 34.50220895241866

$ws.Range("B20").Value = 17.7612965112546
$ws.Range("C20").Value = 7.570408482898304
$ws.Range("D20").Value = 13.89035537652669
$ws.Range("E20").Value = 14.29725412017468
$ws.Range("G20").Value = 3.72777270388296
$ws.Range("J20").Value = 8.664378830694726
$ws.Range("K20").Value = 13.37950851494351
$ws.Range("L20").Value = 12.09245945355128
$ws.Range("M20").Value = 17.79005223252825
$ws.Range("N20").Value = 23.50050592254844
$ws.Range("O20").Value = 34.48190262231593

$ws.Range("B21").Value = 17.95237274089072
$ws.Range("C21").Value = 7.613549240028224
$ws.Range("D21").Value = 13.90488730593504
$ws.Range("E21").Value = 14.27365801355096
$ws.Range("G21").Value = 3.725407559108616
$ws.Range("J21").Value = 8.658585097047764
$ws.Range("K21").Value = 13.51644943563769
$ws.Range("L21").Value = 12.08528226389627
$ws.Range("M21").Value = 17.81760387143082
$ws.Range("N21").Value = 23.42822609178267
$ws.Range("O21").Value = 34.42006915086932

$ws.Range("B22").Value = 18.07895473842606
$ws.Range("C22").Value = 7.641302726027376
$ws.Range("D22").Value = 13.91632480065089
$ws.Range("E22").Value = 14.25954282205873
$ws.Range("G22").Value = 3.723920626836296
$ws.Range("J22").Value = 8.654946730596578
$ws.Range("K22").Value = 13.60716451070469
$ws.Range("L22").Value = 12.08178134197856
$ws.Range("M22").Value = 17.8376967781882
$ws.Range("N22").Value = 23.38257724884481
$ws.Range("O22").Value = 34.3844241068458

$ws.Range("B23").Value = 18.01125087665827
$ws.Range("C23").Value = 7.626531179086794
$ws.Range("D23").Value = 13.9100436849854
$ws.Range("E23").Value = 14.26695701418116
$ws.Range("G23").Value = 3.724708864371489
$ws.Range("J23").Value = 8.656875073778409
$ws.Range("K23").Value = 13.55864457124691
$ws.Range("L23").Value = 12.08354124341047
$ws.Range("M23").Value = 17.82678321865708
$ws.Range("N23").Value = 23.40679574995039
$ws.Range("O23").Value = 34.40301228452126

$ws.Range("B24").Value = 17.75820468772599
$ws.Range("C24").Value = 7.569696053643519
$ws.Range("D24").Value = 13.89015078443405
$ws.Range("E24").Value = 14.29766187548644
$ws.Range("G24").Value = 3.727812357744866
$ws.Range("J24").Value = 8.664476035906429
$ws.Range("K24").Value = 13.37729253001482
$ws.Range("L24").Value = 12.09259690356746
$ws.Range("M24").Value = 17.78963741989874
$ws.Range("N24").Value = 23.50171422529997
$ws.Range("O24").Value = 34.48299370392588

$ws.Range("B25").Value = 17.49367451876977
$ws.Range("C25").Value = 7.50655673190781
$ws.Range("D25").Value = 13.87709677653955
$ws.Range("E25").Value = 14.33640687580432
$ws.Range("G25").Value = 3.731414221187058
$ws.Range("J25").Value = 8.673314793230842
$ws.Range("K25").Value = 13.18766457163126
$ws.Range("L25").Value = 12.10749788894623
$ws.Range("M25").Value = 17.75865091435566
$ws.Range("N25").Value = 23.61096334040537
$ws.Range("O25").Value = 34.58971778258715
